# Actualizando archivo de formatos
# Adds a new "varios" sheet with global material/unit parameters and the
# corresponding workbook-level defined names, then makes it the active sheet.

$wb = $excel.ActiveWorkbook

# --- Create the new "varios" worksheet as the last tab ------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "varios"

# --- Row 1: Young's modulus ----------------------------------------------
$ws.Range("A1").Value = "E"
$ws.Range("B1").Value = 200000000000
$ws.Range("B1").NumberFormat = "0.00E+00"
$ws.Range("C1").Value = "Pa"
$ws.Range("D1").Value = "módulo de Young"

# --- Row 2: Poisson's ratio -----------------------------------------------
$ws.Range("A2").Value = "nu"
$ws.Range("B2").Value = 0.3
$ws.Range("D2").Value = "coeficiente de Poisson"
$ws.Range("G2").Value = "Unidades de fuerza en N"

# --- Row 3: density --------------------------------------------------------
$ws.Range("A3").Value = "rho"
$ws.Range("B3").Formula = "=7850"
$ws.Range("C3").Value = "kg/m³"
$ws.Range("D3").Value = "densidad"
$ws.Range("G3").Value = "Unidades de longitud en m"

# --- Row 4: gravity ----------------------------------------------------------
$ws.Range("A4").Value = "g"
$ws.Range("B4").Value = 9.81
$ws.Range("C4").Value = "m/s²"
$ws.Range("D4").Value = "aceleracion de la gravedad"

# --- Row 5: thickness --------------------------------------------------------
$ws.Range("A5").Value = "espesor"
$ws.Range("B5").Value = 0.01
$ws.Range("C5").Value = "m"

# --- Row 6: length unit label ------------------------------------------------
$ws.Range("A6").Value = "U_LONG"
$ws.Range("B6").Value = "m"

# --- Row 7: force unit label --------------------------------------------------
$ws.Range("A7").Value = "U_FUERZA"
$ws.Range("B7").Value = "N"

# --- Row 8: stress unit label -------------------------------------------------
$ws.Range("A8").Value = "U_ESFUER"
$ws.Range("B8").Value = "Pa"

# --- Row 9: scale factor -------------------------------------------------------
$ws.Range("A9").Value = "ESC_UV"
$ws.Range("B9").Value = 10000

# --- Workbook-level defined names pointing at the new sheet -------------------
$wb.Names.Add('Young', '=varios!$B$1')
$wb.Names.Add('Poisson', '=varios!$B$2')
$wb.Names.Add('rho', '=varios!$B$3')
$wb.Names.Add('g', '=varios!$B$4')
$wb.Names.Add('espesor', '=varios!$B$5')
$wb.Names.Add('U_LONG', '=varios!$B$6')
$wb.Names.Add('U_FUERZA', '=varios!$B$7')
$wb.Names.Add('U_ESFUERZO', '=varios!$B$8')

# --- Refresh the on-screen selection on the other sheets (A1:I19) -------------
foreach ($name in @("xnod", "LaG_mat", "restric", "carga_punt")) {
    [void]$wb.Worksheets.Item($name).Range("A1:I19").Select()
}
[void]$ws.Range("A1:I19").Select()

# --- Make the new sheet the active tab (matches activeTab="7") ----------------
$ws.Activate()
